$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 541, shifting existing rows 541..652 down to 542..653.
$ws.Rows.Item(541).Insert()

# Populate the new row 541 with the new record's data.
$ws.Range("A541").Value = 4
$ws.Range("B541").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C541").Value = "Los Lagos"
$ws.Range("D541").Value = 45211
$ws.Range("D541").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E541").Value = 10
$ws.Range("F541").Value = 100114013
$ws.Range("G541").Value = "Zanahoria"
$ws.Range("H541").Value = "Sin especificar"
$ws.Range("I541").Value = "Primera"
$ws.Range("J541").Value = 300
$ws.Range("K541").Value = 8500
$ws.Range("L541").Value = 8500
$ws.Range("M541").Value = 8500
$ws.Range("N541").Value = "`$/saco 20 kilos"
$ws.Range("O541").Value = "Provincia de Llanquihue"
$ws.Range("P541").Value = 425
$ws.Range("Q541").Value = 20
$ws.Range("R541").Value = "Hortaliza"
